# Update recalculated result values on the "Calc" and "Results" sheets
# (rows 4, 8, 10, 12, 14, 16) to match the refreshed analysis numbers.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Calc")
$ws.Range("AP4").Value = 29.1141
$ws.Range("AQ4").Value = 7.70222028687061
$ws.Range("AW4").Value = 29.6897705849
$ws.Range("AY4").Value = 7.854515623411029
$ws.Range("BC4").Value = 4.368823903971476
$ws.Range("AP8").Value = 16.1988
$ws.Range("AQ8").Value = 4.13341192501129
$ws.Range("AW8").Value = 15.7358106888
$ws.Range("AY8").Value = 4.015269918591311
$ws.Range("BC8").Value = 2.241911907494027
$ws.Range("AP10").Value = 15.9783
$ws.Range("AQ10").Value = 4.181647262770953
$ws.Range("AW10").Value = 15.564272782
$ws.Range("AY10").Value = 4.073314391340145
$ws.Range("BC10").Value = 1.865866594946743
$ws.Range("AP12").Value = 11.3726
$ws.Range("AQ12").Value = 2.830561784105636
$ws.Range("AW12").Value = 11.9606542028
$ws.Range("AY12").Value = 2.977353487320802
$ws.Range("BC12").Value = 1.938401912866515
$ws.Range("AP14").Value = 18.4266
$ws.Range("AQ14").Value = 4.440100712258023
$ws.Range("AW14").Value = 18.264294678
$ws.Range("AY14").Value = 4.401005168639826
$ws.Range("BC14").Value = 2.219418562005929
$ws.Range("AP16").Value = 13.9275
$ws.Range("AQ16").Value = 3.226306367411839
$ws.Range("AW16").Value = 14.2688349709
$ws.Range("AY16").Value = 3.305686005751946
$ws.Range("BC16").Value = 1.800842916945104
$ws = $wb.Worksheets.Item("Results")
$ws.Range("N4").Value = 29.1141
$ws.Range("P4").Value = 29.6897705849
$ws.Range("R4").Value = 4.368823903971476
$ws.Range("N8").Value = 16.1988
$ws.Range("P8").Value = 15.7358106888
$ws.Range("R8").Value = 2.241911907494027
$ws.Range("N10").Value = 15.9783
$ws.Range("P10").Value = 15.564272782
$ws.Range("R10").Value = 1.865866594946743
$ws.Range("N12").Value = 11.3726
$ws.Range("P12").Value = 11.9606542028
$ws.Range("R12").Value = 1.938401912866515
$ws.Range("N14").Value = 18.4266
$ws.Range("P14").Value = 18.264294678
$ws.Range("R14").Value = 2.219418562005929
$ws.Range("N16").Value = 13.9275
$ws.Range("P16").Value = 14.2688349709
$ws.Range("R16").Value = 1.800842916945104
